# Add "Save" column (H) to the worksheet, with header + per-row values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell, matching the style of the other header cells (B1:G1 use style index 1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Per-row "Save" values (0/1 flags)
$values = @(0, 0, 1, 1, 1, 1, 0, 0, 1, 1)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
